$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / changed date) for rows 2-6 from 2023-09-01 to 2023-09-05
$newDate = Get-Date -Year 2023 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2:C6").Value = $newDate.Date
